$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A1: corrected product code from the scraper ---
# Apply bold + wrapped text *before* writing the value, and force text
# formatting so the all-digit code is stored as text (not auto-converted
# to a number), then drop back to General once the text type has stuck.
$ws.Range("A1").NumberFormat = "@"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").WrapText = $true
$ws.Range("A1").Value = "3273114"
$ws.Range("A1").NumberFormat = "General"

# --- A2 (new row): second product code returned for this basket ---
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").WrapText = $true
$ws.Range("A2").Value = "6SL32105BE211UV0"
$ws.Range("A2").NumberFormat = "General"

# --- B1: quantity missing this run (scrapper_Service error handling) ---
# keep the cell present (and touched/formatted) but blank.
$ws.Range("B1").Style = "Normal"
$ws.Range("B1").Value = ""

$ws.Rows.Item(1).RowHeight = 15.65
$ws.Rows.Item(2).RowHeight = 44
$ws.Columns.Item(1).ColumnWidth = 8.56
